$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep column D formatted as text so values containing multiple
# "." separators (e.g. "27.762.42") are not reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "27.762.42"
$ws.Range("E2").Value = "  +1.46%  "
$ws.Range("D3").Value = "1.866.20"
$ws.Range("E3").Value = "  +1.42%  "
$ws.Range("E4").Value = "  +1.28%  "
$ws.Range("D5").Value = "323.34"
$ws.Range("E5").Value = "  +1.67%  "
$ws.Range("D6").Value = "1.032"
$ws.Range("E6").Value = "  +1.11%  "
$ws.Range("D7").Value = "0.4423"
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("D8").Value = "0.3803"
$ws.Range("E8").Value = "  +2.20%  "
$ws.Range("D9").Value = "0.07467"
$ws.Range("E9").Value = "  +1.76%  "
$ws.Range("D10").Value = "0.8877"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").Value = "21.76"
$ws.Range("E11").Value = "  +1.76%  "
$ws.Range("D12").Value = "1.871.54"
$ws.Range("E12").Value = "  -7.80%  "
$ws.Range("D13").Value = "5.552"
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("D14").Value = "6.780"
$ws.Range("E14").Value = "  +1.64%  "
$ws.Range("D15").Value = "0.07212"
$ws.Range("E15").Value = "  +1.32%  "
$ws.Range("D16").Value = "84.39"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").Value = "1.038"
$ws.Range("E17").Value = "  +1.07%  "
$ws.Range("D18").Value = "0.000009137"
$ws.Range("E18").Value = "  +1.53%  "
$ws.Range("D19").Value = "1.032"
$ws.Range("E19").Value = "  +1.01%  "
$ws.Range("D20").Value = "15.58"
$ws.Range("E20").Value = "  +1.05%  "
$ws.Range("D21").Value = "27.768.80"
$ws.Range("E21").Value = "  +1.33%  "
$ws.Range("D22").Value = "5.312"
$ws.Range("E22").Value = "  +1.36%  "
$ws.Range("E23").Value = "  +1.76%  "
$ws.Range("D24").Value = "2.095.54"
$ws.Range("E24").Value = "  -5.13%  "
$ws.Range("D25").Value = "2.021"
$ws.Range("E25").Value = "  +6.09%  "
$ws.Range("D26").Value = "158.68"
$ws.Range("E26").Value = "  +1.27%  "
$ws.Range("E27").Value = "  +1.65%  "
$ws.Range("D28").Value = "1.994"
$ws.Range("E28").Value = "  +3.46%  "
$ws.Range("D29").Value = "5.345"
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Value = "119.25"
$ws.Range("E30").Value = "  +3.26%  "
$ws.Range("D31").Value = "0.09044"
$ws.Range("E31").Value = "  +0.30%  "
$ws.Range("D32").Value = "1.231"
$ws.Range("E32").Value = "  +2.44%  "
$ws.Range("D33").Value = "0.7798"
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "3.032"
$ws.Range("E34").Value = "  +7.01%  "
$ws.Range("D35").Value = "4.593"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("D36").Value = "1.034"
$ws.Range("E36").Value = "  +0.90%  "
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "0.01989"
$ws.Range("E38").Value = "  +1.79%  "
$ws.Range("D39").Value = "0.05359"
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("D40").Value = "2.882"
$ws.Range("E40").Value = "  +3.05%  "
$ws.Range("E41").Value = "  +1.35%  "
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").Value = "6.898"
$ws.Range("E43").Value = "  +5.59%  "
$ws.Range("D44").Value = "8.708"
$ws.Range("E44").Value = "  +2.78%  "
$ws.Range("D45").Value = "110.72"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("D46").Value = "0.06743"
$ws.Range("E46").Value = "  +7.37%  "
$ws.Range("D47").Value = "10.69"
$ws.Range("E47").Value = "  +1.62%  "
$ws.Range("D48").Value = "1.716"
$ws.Range("E48").Value = "  +2.75%  "
$ws.Range("D49").Value = "0.4731"
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Value = "1.921"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").Value = "39.73"
$ws.Range("E51").Value = "  +1.01%  "
